$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New variable-name labels for rows 21-63 (B/C values in these rows stay unchanged;
# only the column-A variable names shift up by one, with "deprawsc" newly inserted).
$names = @(
    "dep_secret",
    "deprawsc",
    "discrim",
    "divers",
    "drug_mar",
    "drugs_yn",
    "dx_adhd",
    "dx_bi",
    "dx_dep",
    "dx_pers",
    "dx_tr",
    "ed_any",
    "env_mh",
    "fincur",
    "finpast",
    "flourish",
    "gad7_impa",
    "gender_noncis",
    "gpa_sr",
    "inf",
    "ins_cover",
    "international",
    "meds_anx",
    "meds_count",
    "meds_dep",
    "meds_mood",
    "meds_sle",
    "meds_sti",
    "military",
    "percneed_cur",
    "persist",
    "psyhx",
    "race",
    "religios",
    "residenc",
    "satisfied_overall",
    "school2_type",
    "sexual",
    "sib_freq",
    "stig_pcv_2",
    "stig_pcv_3",
    "talk",
    "ther_vis"
)

$startRow = 21
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $names[$i]
}

# Remove the now-obsolete last row (64) entirely, shrinking the used range to A1:C63
$ws.Rows.Item(64).Delete() | Out-Null
